# Update the "Table 2.1" baseline-characteristics table with the new PA4
# dataframe's values (larger cohort: N = 87,276 instead of N = 33,143).
# Each statistic in the table is updated via an exact literal Find/Replace
# over the whole document content.
$d = $word.ActiveDocument

$d.Content.Find.Execute("N = 33,143", $true, $false, $false, $false, $false, $true, 1, $false, "N = 87,276", 2) | Out-Null
$d.Content.Find.Execute("96.5 (89.8, 102.8)", $true, $false, $false, $false, $false, $true, 1, $false, "96.5 (89.8, 102.7)", 2) | Out-Null
$d.Content.Find.Execute("164.6 (156.1, 172.7)", $true, $false, $false, $false, $false, $true, 1, $false, "164.7 (156.1, 172.9)", 2) | Out-Null
$d.Content.Find.Execute("465.2 (321.9, 637.8)", $true, $false, $false, $false, $false, $true, 1, $false, "462.2 (318.4, 638.1)", 2) | Out-Null
$d.Content.Find.Execute("235.1 (116.5, 407.3)", $true, $false, $false, $false, $false, $true, 1, $false, "232.5 (114.0, 401.1)", 2) | Out-Null
$d.Content.Find.Execute("587.5 (375.1, 853.5)", $true, $false, $false, $false, $false, $true, 1, $false, "584.2 (371.7, 849.1)", 2) | Out-Null
$d.Content.Find.Execute("738 (2.2)", $true, $false, $false, $false, $false, $true, 1, $false, "1,917 (2.2)", 2) | Out-Null
$d.Content.Find.Execute("470 (1.4)", $true, $false, $false, $false, $false, $true, 1, $false, "1,227 (1.4)", 2) | Out-Null
$d.Content.Find.Execute("63.6 (56.5, 68.6)", $true, $false, $false, $false, $false, $true, 1, $false, "63.2 (56.1, 68.4)", 2) | Out-Null
$d.Content.Find.Execute("994 (3.0)", $true, $false, $false, $false, $false, $true, 1, $false, "2,590 (3.0)", 2) | Out-Null
$d.Content.Find.Execute("32,149 (97)", $true, $false, $false, $false, $false, $true, 1, $false, "84,686 (97)", 2) | Out-Null
$d.Content.Find.Execute("19,169 (58)", $true, $false, $false, $false, $false, $true, 1, $false, "50,371 (58)", 2) | Out-Null
$d.Content.Find.Execute("13,974 (42)", $true, $false, $false, $false, $false, $true, 1, $false, "36,905 (42)", 2) | Out-Null
$d.Content.Find.Execute("2,564 (7.7)", $true, $false, $false, $false, $false, $true, 1, $false, "6,683 (7.7)", 2) | Out-Null
$d.Content.Find.Execute("7,956 (24)", $true, $false, $false, $false, $false, $true, 1, $false, "21,174 (24)", 2) | Out-Null
$d.Content.Find.Execute("7,741 (23)", $true, $false, $false, $false, $false, $true, 1, $false, "20,530 (24)", 2) | Out-Null
$d.Content.Find.Execute("14,882 (45)", $true, $false, $false, $false, $false, $true, 1, $false, "38,889 (45)", 2) | Out-Null
$d.Content.Find.Execute("4,244 (13)", $true, $false, $false, $false, $false, $true, 1, $false, "11,034 (13)", 2) | Out-Null
$d.Content.Find.Execute("7,236 (22)", $true, $false, $false, $false, $false, $true, 1, $false, "18,743 (21)", 2) | Out-Null
$d.Content.Find.Execute("8,680 (26)", $true, $false, $false, $false, $false, $true, 1, $false, "22,949 (26)", 2) | Out-Null
$d.Content.Find.Execute("7,674 (23)", $true, $false, $false, $false, $false, $true, 1, $false, "20,512 (24)", 2) | Out-Null
$d.Content.Find.Execute("2,221 (6.7)", $true, $false, $false, $false, $false, $true, 1, $false, "6,000 (6.9)", 2) | Out-Null
$d.Content.Find.Execute("3,088 (9.3)", $true, $false, $false, $false, $false, $true, 1, $false, "8,038 (9.2)", 2) | Out-Null
$d.Content.Find.Execute("-2.5 (-3.9, -0.3)", $true, $false, $false, $false, $false, $true, 1, $false, "-2.5 (-3.8, -0.2)", 2) | Out-Null
$d.Content.Find.Execute("19,136 (58)", $true, $false, $false, $false, $false, $true, 1, $false, "50,562 (58)", 2) | Out-Null
$d.Content.Find.Execute("11,801 (36)", $true, $false, $false, $false, $false, $true, 1, $false, "30,886 (35)", 2) | Out-Null
$d.Content.Find.Execute("2,206 (6.7)", $true, $false, $false, $false, $false, $true, 1, $false, "5,828 (6.7)", 2) | Out-Null
$d.Content.Find.Execute("1,819 (5.5)", $true, $false, $false, $false, $false, $true, 1, $false, "4,765 (5.5)", 2) | Out-Null
$d.Content.Find.Execute("6,596 (20)", $true, $false, $false, $false, $false, $true, 1, $false, "17,618 (20)", 2) | Out-Null
$d.Content.Find.Execute("8,187 (25)", $true, $false, $false, $false, $false, $true, 1, $false, "21,962 (25)", 2) | Out-Null
$d.Content.Find.Execute("8,744 (26)", $true, $false, $false, $false, $false, $true, 1, $false, "22,953 (26)", 2) | Out-Null
$d.Content.Find.Execute("7,797 (24)", $true, $false, $false, $false, $false, $true, 1, $false, "19,978 (23)", 2) | Out-Null
$d.Content.Find.Execute("23,637 (71)", $true, $false, $false, $false, $false, $true, 1, $false, "62,223 (71)", 2) | Out-Null
$d.Content.Find.Execute("8,369 (25)", $true, $false, $false, $false, $false, $true, 1, $false, "21,985 (25)", 2) | Out-Null
$d.Content.Find.Execute("1,137 (3.4)", $true, $false, $false, $false, $false, $true, 1, $false, "3,068 (3.5)", 2) | Out-Null
$d.Content.Find.Execute("5,976 (18)", $true, $false, $false, $false, $false, $true, 1, $false, "15,701 (18)", 2) | Out-Null
$d.Content.Find.Execute("11,299 (34)", $true, $false, $false, $false, $false, $true, 1, $false, "29,625 (34)", 2) | Out-Null
$d.Content.Find.Execute("15,868 (48)", $true, $false, $false, $false, $false, $true, 1, $false, "41,950 (48)", 2) | Out-Null
$d.Content.Find.Execute("27,412 (83)", $true, $false, $false, $false, $false, $true, 1, $false, "72,160 (83)", 2) | Out-Null
$d.Content.Find.Execute("5,477 (17)", $true, $false, $false, $false, $false, $true, 1, $false, "14,426 (17)", 2) | Out-Null
$d.Content.Find.Execute("254 (0.8)", $true, $false, $false, $false, $false, $true, 1, $false, "690 (0.8)", 2) | Out-Null
$d.Content.Find.Execute("28,209 (85)", $true, $false, $false, $false, $false, $true, 1, $false, "74,100 (85)", 2) | Out-Null
$d.Content.Find.Execute("4,786 (14)", $true, $false, $false, $false, $false, $true, 1, $false, "12,770 (15)", 2) | Out-Null
$d.Content.Find.Execute("148 (0.4)", $true, $false, $false, $false, $false, $true, 1, $false, "406 (0.5)", 2) | Out-Null
$d.Content.Find.Execute("10,769 (32)", $true, $false, $false, $false, $false, $true, 1, $false, "28,516 (33)", 2) | Out-Null
$d.Content.Find.Execute("19,946 (60)", $true, $false, $false, $false, $false, $true, 1, $false, "52,300 (60)", 2) | Out-Null
$d.Content.Find.Execute("2,428 (7.3)", $true, $false, $false, $false, $false, $true, 1, $false, "6,460 (7.4)", 2) | Out-Null
$d.Content.Find.Execute("5,594 (17)", $true, $false, $false, $false, $false, $true, 1, $false, "14,915 (17)", 2) | Out-Null
$d.Content.Find.Execute("24,492 (74)", $true, $false, $false, $false, $false, $true, 1, $false, "64,378 (74)", 2) | Out-Null
$d.Content.Find.Execute("3,057 (9.2)", $true, $false, $false, $false, $false, $true, 1, $false, "7,983 (9.1)", 2) | Out-Null
